# Update cells with recomputed values based on new TPM input data
# (NATMI ligand-receptor edge statistics for Lamc2 -> Itgb1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.1133093333333333
$ws.Range("H2").Value = 0.339928
$ws.Range("I2").Value = 0.02456654176752224
$ws.Range("J2").Value = 0.02456654176752224
$ws.Range("M2").Value = 61.04160633333334
$ws.Range("N2").Value = 183.124819
$ws.Range("O2").Value = 0.2043613460574534
$ws.Range("P2").Value = 0.2043613460574534
$ws.Range("Q2").Value = 6.916583719225779
$ws.Range("R2").Value = 62.24925347303201
$ws.Range("S2").Value = 0.005020451543587496
$ws.Range("T2").Value = 0.005020451543587496

# Row 3
$ws.Range("G3").Value = 0.1133093333333333
$ws.Range("H3").Value = 0.339928
$ws.Range("I3").Value = 0.02456654176752224
$ws.Range("J3").Value = 0.02456654176752224
$ws.Range("O3").Value = 0.3559304658284363
$ws.Range("P3").Value = 0.3559304658284363
$ws.Range("Q3").Value = 12.04642126614933
$ws.Range("R3").Value = 108.417791395344
$ws.Range("S3").Value = 0.008743980655107926
$ws.Range("T3").Value = 0.008743980655107928

# Row 4
$ws.Range("G4").Value = 0.1133093333333333
$ws.Range("H4").Value = 0.339928
$ws.Range("I4").Value = 0.02456654176752224
$ws.Range("J4").Value = 0.02456654176752224
$ws.Range("M4").Value = 131.3384093333333
$ws.Range("N4").Value = 394.015228
$ws.Range("O4").Value = 0.4397081881141102
$ws.Range("P4").Value = 0.4397081881141103
$ws.Range("Q4").Value = 14.88186760262044
$ws.Range("R4").Value = 133.936808423584
$ws.Range("S4").Value = 0.01080210956882681
$ws.Range("T4").Value = 0.01080210956882681

# Row 5
$ws.Range("I5").Value = 0.8380577451911468
$ws.Range("J5").Value = 0.8380577451911468
$ws.Range("M5").Value = 61.04160633333334
$ws.Range("N5").Value = 183.124819
$ws.Range("O5").Value = 0.2043613460574534
$ws.Range("P5").Value = 0.2043613460574534
$ws.Range("Q5").Value = 235.9508558841321
$ws.Range("R5").Value = 2123.557702957189
$ws.Range("S5").Value = 0.1712666088811371
$ws.Range("T5").Value = 0.1712666088811371

# Row 6
$ws.Range("I6").Value = 0.8380577451911468
$ws.Range("J6").Value = 0.8380577451911468
$ws.Range("O6").Value = 0.3559304658284363
$ws.Range("P6").Value = 0.3559304658284363
$ws.Range("S6").Value = 0.2982902836370139
$ws.Range("T6").Value = 0.2982902836370139

# Row 7
$ws.Range("I7").Value = 0.8380577451911468
$ws.Range("J7").Value = 0.8380577451911468
$ws.Range("M7").Value = 131.3384093333333
$ws.Range("N7").Value = 394.015228
$ws.Range("O7").Value = 0.4397081881141102
$ws.Range("P7").Value = 0.4397081881141103
$ws.Range("Q7").Value = 507.6768446006297
$ws.Range("R7").Value = 4569.091601405667
$ws.Range("S7").Value = 0.3685008526729959
$ws.Range("T7").Value = 0.3685008526729959

# Row 8
$ws.Range("G8").Value = 0.6336240000000001
$ws.Range("H8").Value = 1.900872
$ws.Range("I8").Value = 0.1373757130413309
$ws.Range("J8").Value = 0.1373757130413309
$ws.Range("M8").Value = 61.04160633333334
$ws.Range("N8").Value = 183.124819
$ws.Range("O8").Value = 0.2043613460574534
$ws.Range("P8").Value = 0.2043613460574534
$ws.Range("Q8").Value = 38.67742677135201
$ws.Range("R8").Value = 348.0968409421681
$ws.Range("S8").Value = 0.02807428563272885
$ws.Range("T8").Value = 0.02807428563272885

# Row 9
$ws.Range("G9").Value = 0.6336240000000001
$ws.Range("H9").Value = 1.900872
$ws.Range("I9").Value = 0.1373757130413309
$ws.Range("J9").Value = 0.1373757130413309
$ws.Range("O9").Value = 0.3559304658284363
$ws.Range("P9").Value = 0.3559304658284363
$ws.Range("Q9").Value = 67.36339720478401
$ws.Range("R9").Value = 606.270574843056
$ws.Range("S9").Value = 0.04889620153631451
$ws.Range("T9").Value = 0.04889620153631451

# Row 10
$ws.Range("G10").Value = 0.6336240000000001
$ws.Range("H10").Value = 1.900872
$ws.Range("I10").Value = 0.1373757130413309
$ws.Range("J10").Value = 0.1373757130413309
$ws.Range("M10").Value = 131.3384093333333
$ws.Range("N10").Value = 394.015228
$ws.Range("O10").Value = 0.4397081881141102
$ws.Range("P10").Value = 0.4397081881141103
$ws.Range("Q10").Value = 83.219168275424
$ws.Range("R10").Value = 748.972514478816
$ws.Range("S10").Value = 0.06040522587228756
$ws.Range("T10").Value = 0.06040522587228757
